$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for B:E (G is the recalculated sum of B:E for each row)
$data = @{
    2 = @{ B = 3.230985683306322;  C = 1.667794583268128;  D = 3.900430680208489;   E = 0.496779210170732 }
    3 = @{ B = 3.230985683306322;  C = 1.667794583268128;  D = 0.1575252929769615;  E = 0.496779210170732 }
    4 = @{ B = 1.459612070389937;  C = 1.667794583268128;  D = 0.1575252929769615;  E = 0.496779210170732 }
    5 = @{ B = 0.04763786555579896; C = 1.667794583268128; D = 0.1575252929769615;  E = 0.496779210170732 }
    6 = @{ B = 0.04763786555579896; C = 1.667794583268128; D = 0.1575252929769615;  E = 0.496779210170732 }
    7 = @{ B = 1.459612070389937;  C = 0.3127903958511391; D = 0.8054896365839992;  E = 0.496779210170732 }
    8 = @{ B = 3.230985683306322;  C = 1.667794583268128;  D = 0.1575252929769615;  E = 0.496779210170732 }
    9 = @{ B = 3.230985683306322;  C = 1.667794583268128;  D = 0.8054896365839992;  E = 0.496779210170732 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals.B
    $ws.Cells.Item($row, 3).Value = $vals.C
    $ws.Cells.Item($row, 4).Value = $vals.D
    $ws.Cells.Item($row, 5).Value = $vals.E
    $sum = $vals.B + $vals.C + $vals.D + $vals.E
    $ws.Cells.Item($row, 7).Value = $sum
}
